$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 36 data (daily auto push: 2025-09-29 13:37 UTC)
# Column A holds a date-like string ("2025/09/29") that must be stored as
# literal text, not auto-converted into a date serial number. Writing it
# via .Value directly makes the engine "smart match" it into a date (like
# Excel's autocomplete) which also stamps a date NumberFormat style onto
# the cell. To avoid that, write it as a formula that evaluates to the
# text, then convert the formula to its static value via copy / paste-values
# -- this keeps the cell a plain string with no extra formatting applied.
$ws.Range("A36").Formula = '="2025/09/29"'
$ws.Range("A36").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4163) | Out-Null  # xlPasteValues

$ws.Range("B36").Value = "月"
$ws.Range("C36").Value = 20
$ws.Range("D36").Value = 16
